$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": refresh data rows 2-48 (BRVM automated update) ---
$ws1.Cells.Item(2,1).Value = "BRVM - SERVICES PUBLICS"
$ws1.Cells.Item(2,2).Value = 0
$ws1.Cells.Item(2,3).Value = 10
$ws1.Cells.Item(2,4).Value = 4202.29
$ws1.Cells.Item(2,5).Value = 105.88
$ws1.Cells.Item(2,6).Value = "🟡 Observer"
$ws1.Cells.Item(2,7).Value = "➖ Neutre"

$ws1.Cells.Item(3,1).Value = "AIR LIQUIDE CI"
$ws1.Cells.Item(3,2).Value = 0
$ws1.Cells.Item(3,3).Value = 5
$ws1.Cells.Item(3,4).Value = 3480
$ws1.Cells.Item(3,5).Value = 715
$ws1.Cells.Item(3,6).Value = "🟡 Observer"
$ws1.Cells.Item(3,7).Value = "➖ Neutre"

$ws1.Cells.Item(4,1).Value = "NEI-CEDA CI"
$ws1.Cells.Item(4,2).Value = 0
$ws1.Cells.Item(4,3).Value = 5
$ws1.Cells.Item(4,4).Value = 3465
$ws1.Cells.Item(4,5).Value = 725
$ws1.Cells.Item(4,6).Value = "🟡 Observer"
$ws1.Cells.Item(4,7).Value = "➖ Neutre"

$ws1.Cells.Item(5,1).Value = "BRVM - AUTRES SECTEURS"
$ws1.Cells.Item(5,2).Value = 0
$ws1.Cells.Item(5,3).Value = 5
$ws1.Cells.Item(5,4).Value = 3101.62
$ws1.Cells.Item(5,5).Value = 614.21
$ws1.Cells.Item(5,6).Value = "🟡 Observer"
$ws1.Cells.Item(5,7).Value = "➖ Neutre"

$ws1.Cells.Item(6,1).Value = "BRVM - DISTRIBUTION"
$ws1.Cells.Item(6,2).Value = 0
$ws1.Cells.Item(6,3).Value = 5
$ws1.Cells.Item(6,4).Value = 2582.42
$ws1.Cells.Item(6,5).Value = 489.13
$ws1.Cells.Item(6,6).Value = "🟡 Observer"
$ws1.Cells.Item(6,7).Value = "➖ Neutre"

$ws1.Cells.Item(7,1).Value = "BRVM - AGRICULTURE"
$ws1.Cells.Item(7,2).Value = 0
$ws1.Cells.Item(7,3).Value = 5
$ws1.Cells.Item(7,4).Value = 1866.93
$ws1.Cells.Item(7,5).Value = 372.93
$ws1.Cells.Item(7,6).Value = "🟡 Observer"
$ws1.Cells.Item(7,7).Value = "➖ Neutre"

$ws1.Cells.Item(8,1).Value = "BRVM - TRANSPORT"
$ws1.Cells.Item(8,2).Value = 0
$ws1.Cells.Item(8,3).Value = 5
$ws1.Cells.Item(8,4).Value = 1792.96
$ws1.Cells.Item(8,5).Value = 357.37
$ws1.Cells.Item(8,6).Value = "🟡 Observer"
$ws1.Cells.Item(8,7).Value = "➖ Neutre"

$ws1.Cells.Item(9,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(9,2).Value = 0
$ws1.Cells.Item(9,3).Value = 5
$ws1.Cells.Item(9,4).Value = 914.14
$ws1.Cells.Item(9,5).Value = 172.56
$ws1.Cells.Item(9,6).Value = "🟡 Observer"
$ws1.Cells.Item(9,7).Value = "➖ Neutre"

$ws1.Cells.Item(10,1).Value = "BRVM - FINANCES"
$ws1.Cells.Item(10,2).Value = 0
$ws1.Cells.Item(10,3).Value = 5
$ws1.Cells.Item(10,4).Value = 700.33
$ws1.Cells.Item(10,5).Value = 140.3
$ws1.Cells.Item(10,6).Value = "🟡 Observer"
$ws1.Cells.Item(10,7).Value = "➖ Neutre"

$ws1.Cells.Item(11,1).Value = "BRVM-PRESTIGE"
$ws1.Cells.Item(11,2).Value = 0
$ws1.Cells.Item(11,3).Value = 5
$ws1.Cells.Item(11,4).Value = 699.31
$ws1.Cells.Item(11,5).Value = 139.34
$ws1.Cells.Item(11,6).Value = "🟡 Observer"
$ws1.Cells.Item(11,7).Value = "➖ Neutre"

$ws1.Cells.Item(12,1).Value = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(12,2).Value = 0
$ws1.Cells.Item(12,3).Value = 5
$ws1.Cells.Item(12,4).Value = 688.27
$ws1.Cells.Item(12,5).Value = 137.89
$ws1.Cells.Item(12,6).Value = "🟡 Observer"
$ws1.Cells.Item(12,7).Value = "➖ Neutre"

$ws1.Cells.Item(13,1).Value = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(13,2).Value = 0
$ws1.Cells.Item(13,3).Value = 5
$ws1.Cells.Item(13,4).Value = 651.32
$ws1.Cells.Item(13,5).Value = 128.07
$ws1.Cells.Item(13,6).Value = "🟡 Observer"
$ws1.Cells.Item(13,7).Value = "➖ Neutre"

$ws1.Cells.Item(14,1).Value = "BRVM - ENERGIE"
$ws1.Cells.Item(14,2).Value = 0
$ws1.Cells.Item(14,3).Value = 5
$ws1.Cells.Item(14,4).Value = 553.86
$ws1.Cells.Item(14,5).Value = 108.05
$ws1.Cells.Item(14,6).Value = "🟡 Observer"
$ws1.Cells.Item(14,7).Value = "➖ Neutre"

$ws1.Cells.Item(15,1).Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Cells.Item(15,2).Value = 0
$ws1.Cells.Item(15,3).Value = 5
$ws1.Cells.Item(15,4).Value = 476.66
$ws1.Cells.Item(15,5).Value = 95.13
$ws1.Cells.Item(15,6).Value = "🟡 Observer"
$ws1.Cells.Item(15,7).Value = "➖ Neutre"

$ws1.Cells.Item(16,1).Value = "BRVM-PRINCIPAL                    (**)"
$ws1.Cells.Item(16,2).Value = 0
$ws1.Cells.Item(16,3).Value = 2
$ws1.Cells.Item(16,4).Value = 410.25
$ws1.Cells.Item(16,5).Value = 204.33
$ws1.Cells.Item(16,6).Value = "🟡 Observer"
$ws1.Cells.Item(16,7).Value = "➖ Neutre"

$ws1.Cells.Item(17,1).Value = "BRVM - INDUSTRIE                 (**)"
$ws1.Cells.Item(17,2).Value = 0
$ws1.Cells.Item(17,3).Value = 1
$ws1.Cells.Item(17,4).Value = 219.23
$ws1.Cells.Item(17,5).Value = 219.23
$ws1.Cells.Item(17,6).Value = "🟡 Observer"
$ws1.Cells.Item(17,7).Value = "➖ Neutre"

$ws1.Cells.Item(18,1).Value = "BRVM - INDUSTRIE                  (**)"
$ws1.Cells.Item(18,2).Value = 0
$ws1.Cells.Item(18,3).Value = 1
$ws1.Cells.Item(18,4).Value = 218.66
$ws1.Cells.Item(18,5).Value = 218.66
$ws1.Cells.Item(18,6).Value = "🟡 Observer"
$ws1.Cells.Item(18,7).Value = "➖ Neutre"

$ws1.Cells.Item(19,1).Value = "BRVM - INDUSTRIE              (**)"
$ws1.Cells.Item(19,2).Value = 0
$ws1.Cells.Item(19,3).Value = 1
$ws1.Cells.Item(19,4).Value = 214.53
$ws1.Cells.Item(19,5).Value = 214.53
$ws1.Cells.Item(19,6).Value = "🟡 Observer"
$ws1.Cells.Item(19,7).Value = "➖ Neutre"

$ws1.Cells.Item(20,1).Value = "BRVM-PRINCIPAL                (**)"
$ws1.Cells.Item(20,2).Value = 0
$ws1.Cells.Item(20,3).Value = 1
$ws1.Cells.Item(20,4).Value = 205.03
$ws1.Cells.Item(20,5).Value = 205.03
$ws1.Cells.Item(20,6).Value = "🟡 Observer"
$ws1.Cells.Item(20,7).Value = "➖ Neutre"

$ws1.Cells.Item(21,1).Value = "BRVM - CONSOMMATION DE BASE        (**)"
$ws1.Cells.Item(21,2).Value = 0
$ws1.Cells.Item(21,3).Value = 1
$ws1.Cells.Item(21,4).Value = 194.82
$ws1.Cells.Item(21,5).Value = 194.82
$ws1.Cells.Item(21,6).Value = "🟡 Observer"
$ws1.Cells.Item(21,7).Value = "➖ Neutre"

$ws1.Cells.Item(22,1).Value = "BRVM - CONSOMMATION DE BASE          (**)"
$ws1.Cells.Item(22,2).Value = 0
$ws1.Cells.Item(22,3).Value = 1
$ws1.Cells.Item(22,4).Value = 193.9
$ws1.Cells.Item(22,5).Value = 193.9
$ws1.Cells.Item(22,6).Value = "🟡 Observer"
$ws1.Cells.Item(22,7).Value = "➖ Neutre"

$ws1.Cells.Item(23,1).Value = "BRVM - CONSOMMATION DE BASE               (**)"
$ws1.Cells.Item(23,2).Value = 0
$ws1.Cells.Item(23,3).Value = 1
$ws1.Cells.Item(23,4).Value = 191.47
$ws1.Cells.Item(23,5).Value = 191.47
$ws1.Cells.Item(23,6).Value = "🟡 Observer"
$ws1.Cells.Item(23,7).Value = "➖ Neutre"

$ws1.Cells.Item(24,1).Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Cells.Item(24,2).Value = 4
$ws1.Cells.Item(24,3).Value = 0
$ws1.Cells.Item(24,4).Value = 25.86
$ws1.Cells.Item(24,5).Value = 3.75
$ws1.Cells.Item(24,6).Value = "🟢 Achat"
$ws1.Cells.Item(24,7).Value = "✅ Renforcer"

$ws1.Cells.Item(25,1).Value = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(25,2).Value = 3
$ws1.Cells.Item(25,3).Value = 0
$ws1.Cells.Item(25,4).Value = 19.36
$ws1.Cells.Item(25,5).Value = 7.21
$ws1.Cells.Item(25,6).Value = "🟢 Achat"
$ws1.Cells.Item(25,7).Value = "✅ Renforcer"

$ws1.Cells.Item(26,1).Value = "SAFCA CI (SAFC)"
$ws1.Cells.Item(26,2).Value = 3
$ws1.Cells.Item(26,3).Value = 2
$ws1.Cells.Item(26,4).Value = 7.45
$ws1.Cells.Item(26,5).Value = -7.49
$ws1.Cells.Item(26,6).Value = "🟢 Achat"
$ws1.Cells.Item(26,7).Value = "✅ Renforcer"

$ws1.Cells.Item(27,1).Value = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(27,2).Value = 1
$ws1.Cells.Item(27,3).Value = 0
$ws1.Cells.Item(27,4).Value = 7
$ws1.Cells.Item(27,5).Value = 7
$ws1.Cells.Item(27,6).Value = "🟡 Observer"
$ws1.Cells.Item(27,7).Value = "➖ Neutre"

$ws1.Cells.Item(28,1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Cells.Item(28,2).Value = 2
$ws1.Cells.Item(28,3).Value = 0
$ws1.Cells.Item(28,4).Value = 5.35
$ws1.Cells.Item(28,5).Value = 2.29
$ws1.Cells.Item(28,6).Value = "🟡 Observer"
$ws1.Cells.Item(28,7).Value = "➖ Neutre"

$ws1.Cells.Item(29,1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(29,2).Value = 1
$ws1.Cells.Item(29,3).Value = 0
$ws1.Cells.Item(29,4).Value = 5
$ws1.Cells.Item(29,5).Value = 5
$ws1.Cells.Item(29,6).Value = "🟡 Observer"
$ws1.Cells.Item(29,7).Value = "➖ Neutre"

$ws1.Cells.Item(30,1).Value = "SICABLE CI (CABC)"
$ws1.Cells.Item(30,2).Value = 2
$ws1.Cells.Item(30,3).Value = 2
$ws1.Cells.Item(30,4).Value = 4.74
$ws1.Cells.Item(30,5).Value = -7.37
$ws1.Cells.Item(30,6).Value = "🟡 Observer"
$ws1.Cells.Item(30,7).Value = "👀 À surveiller"

$ws1.Cells.Item(31,1).Value = "SOLIBRA CI (SLBC)"
$ws1.Cells.Item(31,2).Value = 1
$ws1.Cells.Item(31,3).Value = 1
$ws1.Cells.Item(31,4).Value = 3.7
$ws1.Cells.Item(31,5).Value = 7.48
$ws1.Cells.Item(31,6).Value = "🟡 Observer"
$ws1.Cells.Item(31,7).Value = "👀 À surveiller"

$ws1.Cells.Item(32,1).Value = "FILTISAC CI (FTSC)"
$ws1.Cells.Item(32,2).Value = 1
$ws1.Cells.Item(32,3).Value = 1
$ws1.Cells.Item(32,4).Value = 3.47
$ws1.Cells.Item(32,5).Value = -3.07
$ws1.Cells.Item(32,6).Value = "🟡 Observer"
$ws1.Cells.Item(32,7).Value = "👀 À surveiller"

$ws1.Cells.Item(33,1).Value = "SAPH CI (SPHC)"
$ws1.Cells.Item(33,2).Value = 1
$ws1.Cells.Item(33,3).Value = 0
$ws1.Cells.Item(33,4).Value = 1.44
$ws1.Cells.Item(33,5).Value = 1.44
$ws1.Cells.Item(33,6).Value = "🟡 Observer"
$ws1.Cells.Item(33,7).Value = "➖ Neutre"

$ws1.Cells.Item(34,1).Value = "SICOR CI (SICC)"
$ws1.Cells.Item(34,2).Value = 1
$ws1.Cells.Item(34,3).Value = 1
$ws1.Cells.Item(34,4).Value = 0.55
$ws1.Cells.Item(34,5).Value = 7.46
$ws1.Cells.Item(34,6).Value = "🟡 Observer"
$ws1.Cells.Item(34,7).Value = "👀 À surveiller"

$ws1.Cells.Item(35,1).Value = "SOGB CI (SOGC)"
$ws1.Cells.Item(35,2).Value = 1
$ws1.Cells.Item(35,3).Value = 1
$ws1.Cells.Item(35,4).Value = 0.49
$ws1.Cells.Item(35,5).Value = 2.07
$ws1.Cells.Item(35,6).Value = "🟡 Observer"
$ws1.Cells.Item(35,7).Value = "👀 À surveiller"

$ws1.Cells.Item(36,1).Value = "TOTAL"
$ws1.Cells.Item(36,2).Value = 0
$ws1.Cells.Item(36,3).Value = 5
$ws1.Cells.Item(36,4).Value = 0
$ws1.Cells.Item(36,5).Value = 0
$ws1.Cells.Item(36,6).Value = "🟡 Observer"
$ws1.Cells.Item(36,7).Value = "➖ Neutre"

$ws1.Cells.Item(37,1).Value = "ONATEL BF (ONTBF)"
$ws1.Cells.Item(37,2).Value = 1
$ws1.Cells.Item(37,3).Value = 1
$ws1.Cells.Item(37,4).Value = -0.65
$ws1.Cells.Item(37,5).Value = 5.49
$ws1.Cells.Item(37,6).Value = "🟡 Observer"
$ws1.Cells.Item(37,7).Value = "👀 À surveiller"

$ws1.Cells.Item(38,1).Value = "NEI-CEDA CI (NEIC)"
$ws1.Cells.Item(38,2).Value = 1
$ws1.Cells.Item(38,3).Value = 1
$ws1.Cells.Item(38,4).Value = -1.45
$ws1.Cells.Item(38,5).Value = 1.45
$ws1.Cells.Item(38,6).Value = "🟡 Observer"
$ws1.Cells.Item(38,7).Value = "👀 À surveiller"

$ws1.Cells.Item(39,1).Value = "BANK OF AFRICA ML (BOAM)"
$ws1.Cells.Item(39,2).Value = 0
$ws1.Cells.Item(39,3).Value = 1
$ws1.Cells.Item(39,4).Value = -2.19
$ws1.Cells.Item(39,5).Value = -2.19
$ws1.Cells.Item(39,6).Value = "🟡 Observer"
$ws1.Cells.Item(39,7).Value = "➖ Neutre"

$ws1.Cells.Item(40,1).Value = "BANK OF AFRICA NG (BOAN)"
$ws1.Cells.Item(40,2).Value = 0
$ws1.Cells.Item(40,3).Value = 1
$ws1.Cells.Item(40,4).Value = -3.7
$ws1.Cells.Item(40,5).Value = -3.7
$ws1.Cells.Item(40,6).Value = "🟡 Observer"
$ws1.Cells.Item(40,7).Value = "➖ Neutre"

$ws1.Cells.Item(41,1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Cells.Item(41,2).Value = 0
$ws1.Cells.Item(41,3).Value = 1
$ws1.Cells.Item(41,4).Value = -3.85
$ws1.Cells.Item(41,5).Value = -3.85
$ws1.Cells.Item(41,6).Value = "🟡 Observer"
$ws1.Cells.Item(41,7).Value = "➖ Neutre"

$ws1.Cells.Item(42,1).Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Cells.Item(42,2).Value = 0
$ws1.Cells.Item(42,3).Value = 1
$ws1.Cells.Item(42,4).Value = -3.91
$ws1.Cells.Item(42,5).Value = -3.91
$ws1.Cells.Item(42,6).Value = "🟡 Observer"
$ws1.Cells.Item(42,7).Value = "➖ Neutre"

$ws1.Cells.Item(43,1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(43,2).Value = 0
$ws1.Cells.Item(43,3).Value = 1
$ws1.Cells.Item(43,4).Value = -5.66
$ws1.Cells.Item(43,5).Value = -5.66
$ws1.Cells.Item(43,6).Value = "🟡 Observer"
$ws1.Cells.Item(43,7).Value = "➖ Neutre"

$ws1.Cells.Item(44,1).Value = "AIR LIQUIDE CI (SIVC)"
$ws1.Cells.Item(44,2).Value = 0
$ws1.Cells.Item(44,3).Value = 2
$ws1.Cells.Item(44,4).Value = -5.74
$ws1.Cells.Item(44,5).Value = -3.5
$ws1.Cells.Item(44,6).Value = "🟡 Observer"
$ws1.Cells.Item(44,7).Value = "➖ Neutre"

$ws1.Cells.Item(45,1).Value = "CFAO MOTORS CI (CFAC)"
$ws1.Cells.Item(45,2).Value = 2
$ws1.Cells.Item(45,3).Value = 3
$ws1.Cells.Item(45,4).Value = -5.94
$ws1.Cells.Item(45,5).Value = -5.75
$ws1.Cells.Item(45,6).Value = "🔴 Vente"
$ws1.Cells.Item(45,7).Value = "⚠️ Risque de décrochage"

$ws1.Cells.Item(46,1).Value = "CIE CI (CIEC)"
$ws1.Cells.Item(46,2).Value = 0
$ws1.Cells.Item(46,3).Value = 1
$ws1.Cells.Item(46,4).Value = -6.3
$ws1.Cells.Item(46,5).Value = -6.3
$ws1.Cells.Item(46,6).Value = "🟡 Observer"
$ws1.Cells.Item(46,7).Value = "➖ Neutre"

$ws1.Cells.Item(47,1).Value = "SETAO CI (STAC)"
$ws1.Cells.Item(47,2).Value = 0
$ws1.Cells.Item(47,3).Value = 2
$ws1.Cells.Item(47,4).Value = -10.33
$ws1.Cells.Item(47,5).Value = -7.08
$ws1.Cells.Item(47,6).Value = "🟡 Observer"
$ws1.Cells.Item(47,7).Value = "➖ Neutre"

$ws1.Cells.Item(48,1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Cells.Item(48,2).Value = 0
$ws1.Cells.Item(48,3).Value = 2
$ws1.Cells.Item(48,4).Value = -14.35
$ws1.Cells.Item(48,5).Value = -7.35
$ws1.Cells.Item(48,6).Value = "🟡 Observer"
$ws1.Cells.Item(48,7).Value = "➖ Neutre"

# --- Sheet "Top_YTD": refresh Progression YTD (%) values rows 2-11 ---
$ws2.Cells.Item(2,2).Value = 157858674.03
$ws2.Cells.Item(3,2).Value = 3190364.2
$ws2.Cells.Item(4,2).Value = 3131428.4
$ws2.Cells.Item(5,2).Value = 1939047.31
$ws2.Cells.Item(6,2).Value = 888944.27
$ws2.Cells.Item(7,2).Value = 237598.05
$ws2.Cells.Item(8,2).Value = 202686.86
$ws2.Cells.Item(9,2).Value = 17968.68
$ws2.Cells.Item(10,2).Value = 7873.36
$ws2.Cells.Item(11,2).Value = 7839.3
